# Update Pooh Points site
# - Column G ("status") values of "End of 2nd Half" become "Final"
# - Column G width shrinks from 17 to 8 (no longer needs the extra room)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "End of 2nd Half") {
        $cell.Value2 = "Final"
    }
}

$ws.Columns.Item(7).ColumnWidth = 7.17
